# Update SBTM Report Template v.1.1.xlsx
# Removes a duplicated/erroneous SBTM session log entry that had leaked
# into the "Summary" sheet's running totals (rows 21-23), clearing the
# now-stale numbers while keeping the surrounding table structure intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- Row 21: clear the stray Sessions/Bugs running totals (keep formatting) ---
$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()

# --- Row 22: C22 is removed outright (content + formatting); D22 keeps its
#     formatting but loses its value ---
$ws.Range("C22").Clear()
$ws.Range("D22").ClearContents()

# --- Row 23: this whole log entry is a duplicate of row 22's data; drop the
#     duplicated fields entirely, but leave K23/N23/O23 behind (still
#     formatted, now empty) since the template keeps those columns primed
#     for the next real entry ---
$ws.Range("H23").Clear()
$ws.Range("I23").Clear()
$ws.Range("J23").Clear()
$ws.Range("L23").Clear()
$ws.Range("M23").Clear()
$ws.Range("K23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("O23").ClearContents()
